$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "35.529.82"
$ws.Cells.Item(2, 5).Value = "  -2.63%  "
$ws.Cells.Item(3, 4).Value = "1.975.21"
$ws.Cells.Item(3, 5).Value = "  -3.94%  "
$ws.Cells.Item(4, 5).Value = "  -0.01%  "
$ws.Cells.Item(5, 4).Value = "'244.90"
$ws.Cells.Item(5, 5).Value = "  +1.11%  "
$ws.Cells.Item(6, 4).Value = "'0.636"
$ws.Cells.Item(6, 5).Value = "  -4.51%  "
$ws.Cells.Item(7, 4).Value = "'57.36"
$ws.Cells.Item(7, 5).Value = "  +4.96%  "
$ws.Cells.Item(8, 5).Value = "  +0.04%  "
$ws.Cells.Item(9, 4).Value = "'58.40"
$ws.Cells.Item(9, 5).Value = "  -0.12%  "
$ws.Cells.Item(10, 5).Value = "  +0.22%  "
$ws.Cells.Item(11, 5).Value = "  -2.53%  "
$ws.Cells.Item(12, 5).Value = "  -3.04%  "
$ws.Cells.Item(13, 4).Value = "'0.946"
$ws.Cells.Item(13, 5).Value = "  +5.99%  "
$ws.Cells.Item(14, 4).Value = "'14.21"
$ws.Cells.Item(14, 5).Value = "  -3.30%  "
$ws.Cells.Item(15, 4).Value = "2.262.47"
$ws.Cells.Item(15, 5).Value = "  -3.96%  "
$ws.Cells.Item(16, 5).Value = "  -1.79%  "
$ws.Cells.Item(17, 4).Value = "1.978.47"
$ws.Cells.Item(17, 5).Value = "  -3.84%  "
$ws.Cells.Item(18, 4).Value = "'17.57"
$ws.Cells.Item(18, 5).Value = "  +5.02%  "
$ws.Cells.Item(19, 4).Value = "35.400.76"
$ws.Cells.Item(19, 5).Value = "  -2.77%  "
$ws.Cells.Item(20, 4).Value = "'71.41"
$ws.Cells.Item(21, 4).Value = "0.0$([char]0x2083)0838"
$ws.Cells.Item(21, 5).Value = "  -1.91%  "
$ws.Cells.Item(22, 4).Value = "'232.24"
$ws.Cells.Item(22, 5).Value = "  -2.62%  "
$ws.Cells.Item(23, 5).Value = "  -2.23%  "
$ws.Cells.Item(24, 5).Value = "  +0.08%  "
$ws.Cells.Item(25, 5).Value = "  +20.73%  "
$ws.Cells.Item(26, 5).Value = "  -1.77%  "
$ws.Cells.Item(27, 4).Value = "'163.51"
$ws.Cells.Item(27, 5).Value = "  +0.47%  "
$ws.Cells.Item(28, 4).Value = "'9.02"
$ws.Cells.Item(28, 5).Value = "  -3.64%  "
$ws.Cells.Item(29, 4).Value = "'19.18"
$ws.Cells.Item(29, 5).Value = "  -4.84%  "
$ws.Cells.Item(30, 5).Value = "  -2.64%  "
$ws.Cells.Item(31, 4).Value = "'4.86"
$ws.Cells.Item(31, 5).Value = "  -3.51%  "
$ws.Cells.Item(32, 4).Value = "'1.12"
$ws.Cells.Item(32, 5).Value = "  -4.70%  "
$ws.Cells.Item(33, 5).Value = "  -0.49%  "
$ws.Cells.Item(34, 4).Value = "'0.0910"
$ws.Cells.Item(34, 5).Value = "  +10.45%  "
$ws.Cells.Item(35, 4).Value = "'4.26"
$ws.Cells.Item(35, 5).Value = "  -5.16%  "
$ws.Cells.Item(36, 5).Value = "  +7.76%  "
$ws.Cells.Item(37, 5).Value = "  +0.05%  "
$ws.Cells.Item(38, 4).Value = "'1.77"
$ws.Cells.Item(38, 5).Value = "  -4.54%  "
$ws.Cells.Item(39, 5).Value = "  +4.96%  "
$ws.Cells.Item(40, 5).Value = "  -3.04%  "
$ws.Cells.Item(41, 5).Value = "  +0.97%  "
$ws.Cells.Item(42, 5).Value = "  -1.53%  "
$ws.Cells.Item(43, 5).Value = "  -2.51%  "
$ws.Cells.Item(44, 4).Value = "'91.18"
$ws.Cells.Item(44, 5).Value = "  -2.54%  "
$ws.Cells.Item(45, 4).Value = "1.378.57"
$ws.Cells.Item(45, 5).Value = "  -0.79%  "
$ws.Cells.Item(46, 4).Value = "'15.85"
$ws.Cells.Item(46, 5).Value = "  +1.20%  "
$ws.Cells.Item(47, 4).Value = "'0.0882"
$ws.Cells.Item(47, 5).Value = "  -2.10%  "
$ws.Cells.Item(48, 4).Value = "'7.47"
$ws.Cells.Item(48, 5).Value = "  +2.28%  "
$ws.Cells.Item(49, 5).Value = "  +1.07%  "
$ws.Cells.Item(50, 2).Value = "MultiversX"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Cells.Item(50, 4).Value = "'45.65"
$ws.Cells.Item(50, 5).Value = "  +0.74%  "
$ws.Cells.Item(51, 4).Value = "'2.24"
$ws.Cells.Item(51, 5).Value = "  -0.75%  "

# Reset quote-prefix styling picked up from the text-forcing trick above
# so the edited cells keep the workbook's original (unstyled) appearance.
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(51, 4).Style = "Normal"
